$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '51.941.92'
$ws.Cells.Item(2, 5).Value = '  -14.32%  '

$ws.Cells.Item(3, 4).Value = '2.279.01'
$ws.Cells.Item(3, 5).Value = '  -21.54%  '

$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 4).Value = '441.15'
$ws.Cells.Item(5, 5).Value = '  -16.21%  '

$ws.Cells.Item(6, 4).Value = '118.80'
$ws.Cells.Item(6, 5).Value = '  -17.07%  '

$ws.Cells.Item(7, 4).Value = '0.998'
$ws.Cells.Item(7, 5).Value = '  -0.06%  '

$ws.Cells.Item(8, 4).Value = '0.460'
$ws.Cells.Item(8, 5).Value = '  -16.05%  '

$ws.Cells.Item(9, 4).Value = '2.268.33'
$ws.Cells.Item(9, 5).Value = '  -22.09%  '

$ws.Cells.Item(10, 4).Value = '5.23'
$ws.Cells.Item(10, 5).Value = '  -12.56%  '

$ws.Cells.Item(11, 4).Value = '0.0859'
$ws.Cells.Item(11, 5).Value = '  -19.88%  '

$ws.Cells.Item(12, 4).Value = '0.299'
$ws.Cells.Item(12, 5).Value = '  -16.77%  '

$ws.Cells.Item(13, 5).Value = '  -6.60%  '

$ws.Cells.Item(14, 4).Value = '51.915.81'
$ws.Cells.Item(14, 5).Value = '  -14.31%  '

$ws.Cells.Item(15, 4).Value = '18.53'
$ws.Cells.Item(15, 5).Value = '  -17.91%  '

$ws.Cells.Item(16, 4).Value = '0.0000116'
$ws.Cells.Item(16, 5).Value = '  -17.67%  '

$ws.Cells.Item(17, 4).Value = '2.294.13'
$ws.Cells.Item(17, 5).Value = '  -21.14%  '

$ws.Cells.Item(18, 4).Value = '3.92'
$ws.Cells.Item(18, 5).Value = '  -21.38%  '

$ws.Cells.Item(19, 4).Value = '294.82'
$ws.Cells.Item(19, 5).Value = '  -15.77%  '

$ws.Cells.Item(20, 2).Value = 'Dai'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(20, 4).Value = '0.999'
$ws.Cells.Item(20, 5).Value = '  -0.19%  '

$ws.Cells.Item(21, 2).Value = 'Chainlink'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(21, 4).Value = '8.72'
$ws.Cells.Item(21, 5).Value = '  -24.77%  '

$ws.Cells.Item(22, 4).Value = '5.61'

$ws.Cells.Item(23, 4).Value = '5.08'
$ws.Cells.Item(23, 5).Value = '  -22.36%  '

$ws.Cells.Item(24, 4).Value = '52.93'
$ws.Cells.Item(24, 5).Value = '  -18.14%  '

$ws.Cells.Item(25, 4).Value = '0.362'
$ws.Cells.Item(25, 5).Value = '  -20.03%  '

$ws.Cells.Item(26, 4).Value = '0.142'
$ws.Cells.Item(26, 5).Value = '  -20.34%  '

$ws.Cells.Item(27, 2).Value = 'USDe'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(27, 4).Value = '0.999'
$ws.Cells.Item(27, 5).Value = '  -0.07%  '

$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).Value = '6.83'
$ws.Cells.Item(28, 5).Value = '  -12.64%  '

$ws.Cells.Item(29, 4).Value = '0.0₃0647'
$ws.Cells.Item(29, 5).Value = '  -24.07%  '

$ws.Cells.Item(30, 4).Value = '142.09'
$ws.Cells.Item(30, 5).Value = '  -6.23%  '

$ws.Cells.Item(31, 4).Value = '16.58'
$ws.Cells.Item(31, 5).Value = '  -15.31%  '

$ws.Cells.Item(32, 4).Value = '1.31'
$ws.Cells.Item(32, 5).Value = '  -21.41%  '

$ws.Cells.Item(33, 4).Value = '4.64'
$ws.Cells.Item(33, 5).Value = '  -16.77%  '

$ws.Cells.Item(34, 4).Value = '0.807'
$ws.Cells.Item(34, 5).Value = '  -19.16%  '

$ws.Cells.Item(35, 4).Value = '3.34'
$ws.Cells.Item(35, 5).Value = '  -22.45%  '

$ws.Cells.Item(36, 4).Value = '0.991'
$ws.Cells.Item(36, 5).Value = '  -0.66%  '

$ws.Cells.Item(37, 4).Value = '0.979'
$ws.Cells.Item(37, 5).Value = '  -18.13%  '

$ws.Cells.Item(38, 4).Value = '31.62'
$ws.Cells.Item(38, 5).Value = '  -16.07%  '

$ws.Cells.Item(39, 5).Value = '  -1.73%  '

$ws.Cells.Item(40, 2).Value = 'Hedera'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(40, 4).Value = '0.0501'
$ws.Cells.Item(40, 5).Value = '  -13.72%  '

$ws.Cells.Item(41, 2).Value = 'Mantle'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(41, 4).Value = '0.547'
$ws.Cells.Item(41, 5).Value = '  -15.67%  '

$ws.Cells.Item(42, 4).Value = '3.07'
$ws.Cells.Item(42, 5).Value = '  -17.47%  '

$ws.Cells.Item(43, 4).Value = '1.895.53'
$ws.Cells.Item(43, 5).Value = '  -17.30%  '

$ws.Cells.Item(44, 4).Value = '1.15'
$ws.Cells.Item(44, 5).Value = '  -21.74%  '

$ws.Cells.Item(45, 4).Value = '0.0811'
$ws.Cells.Item(45, 5).Value = '  -11.67%  '

$ws.Cells.Item(46, 4).Value = '0.0202'
$ws.Cells.Item(46, 5).Value = '  -14.84%  '

$ws.Cells.Item(47, 4).Value = '4.10'
$ws.Cells.Item(47, 5).Value = '  -17.04%  '

$ws.Cells.Item(48, 4).Value = '15.36'
$ws.Cells.Item(48, 5).Value = '  -24.92%  '

$ws.Cells.Item(49, 5).Value = '  -5.21%  '

$ws.Cells.Item(50, 4).Value = '4.41'
$ws.Cells.Item(50, 5).Value = '  -13.50%  '

$ws.Cells.Item(51, 4).Value = '14.71'
$ws.Cells.Item(51, 5).Value = '  -19.57%  '
